# Work Log.docx edit script
# Adds new work-log entries (Jan 14 - Feb 12, 2014) describing the
# server / game-engine work, per the commit "Added server and game
# engine to repository".
#
# Technique notes:
#  - Directly mutating the text of a run that sits next to another
#    run with identical formatting causes the engine to coalesce the
#    two runs into one (normal Word run-merge behaviour). Some of the
#    target headings must keep their multi-run shape (e.g. "February 12"
#    + ", 2014" as two runs), so those are rewritten with
#    Range.InsertXML against a wrapping <w:p> on the *exact* text
#    span (found via Find), which replaces content without touching
#    the neighbouring run.
#  - Range.InsertXML on a COLLAPSED range splices the given runs in
#    place; if the supplied fragment contains more than one <w:p>, the
#    LAST <w:p> in the fragment is merged into the paragraph that sits
#    at the insertion point (its runs become the leading runs of that
#    paragraph). To insert a clean run of whole new paragraphs without
#    disturbing the following paragraph, the fragment ends with an
#    empty dummy <w:p/> which absorbs that merge, and the now-stray
#    empty paragraph is deleted afterwards.

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-Range([string]$text) {
    $r = $d.Content
    $null = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    # Re-wrap as a fresh Range: Range objects still tied to the
    # whole-story Range (Document.Content) keep "insert at point"
    # semantics for InsertXML even after Find narrows Start/End, so a
    # literal sub-range is required to get "replace this span" instead.
    return $d.Range($r.Start, $r.End)
}

# ---------------------------------------------------------------
# Step 1: drop the bookmark that currently sits on "Started defining
# requirements." -- it will be re-created on the new
# "Implemented protocol for server" paragraph further down. (Text
# replacement alone keeps a bookmark that touches the replaced span,
# so it is removed explicitly via the Bookmarks collection instead.)
# ---------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------
# Step 2: rename the last heading, "January 14" -> "February 12",
# keeping it and the ", 2014" run as two separate runs.
# ---------------------------------------------------------------
$r = Find-Range("January 14")
$xml = "<w:p $wns><w:r><w:t>February 12</w:t></w:r></w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------
# Step 3: rename the last body paragraph, "Continued defining
# requirements." -> "Implemented protocol for game engine" (this is
# still the only occurrence of that sentence at this point).
# ---------------------------------------------------------------
$r = Find-Range("Continued defining requirements.")
$xml = "<w:p $wns><w:r><w:t>Implemented protocol for game engine</w:t></w:r></w:p>"
$r.InsertXML($xml)

# ---------------------------------------------------------------
# Step 4: add a trailing empty paragraph, mirroring the diff's extra
# <w:p/> before the final (pre-existing) empty paragraph / sectPr.
# ---------------------------------------------------------------
$r = Find-Range("Implemented protocol for game engine")
$insertPoint = $d.Range($r.End, $r.End)
$insertPoint.InsertXML("<w:p $wns/>")

# ---------------------------------------------------------------
# Step 5: insert the whole new block of entries (Jan 14 .. Feb 12
# prep material) right after "Started defining requirements.", i.e.
# right before the "February 12, 2014" heading renamed in step 2.
# ---------------------------------------------------------------
$newBlock = @"
<w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>January 14, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Continued defining requirements.</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>January 15, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Continued defining requirements.</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>January 20, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>API design</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>January 27, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Game Design</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>February 5, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Researched Python networking and websockets</w:t></w:r></w:p><w:p $wns><w:r><w:t>Implemented websocket chat server</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>February 7</w:t></w:r><w:r><w:t>, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Installed pypi and Tornado</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">February </w:t></w:r><w:r><w:t>8</w:t></w:r><w:r><w:t>, 2014</w:t></w:r></w:p><w:p $wns><w:r><w:t>Built simple echo server in tornado.  Debugged networking issues in HTTP with Fiddler</w:t></w:r></w:p><w:p $wns><w:r><w:t>Implemented protocol for server</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/></w:p><w:p $wns/>
"@

$r = Find-Range("Started defining requirements.")
$insertPoint = $d.Range($r.End, $r.End)
$insertPoint.InsertXML($newBlock)

# The fragment's trailing dummy <w:p/> merged into the "February 12,
# 2014" heading paragraph, which pushed the heading's own runs after
# it, leaving a stray now-empty paragraph right before it. Find it by
# locating "February 12" again and deleting the empty paragraph that
# directly precedes it.
$heading = Find-Range("February 12")
$headingPara = $heading.Paragraphs.First
$strayIndex = $headingPara.Index - 1
$stray = $d.Paragraphs($strayIndex)
if ($stray.Range.Text.Trim().Length -eq 0) {
    $stray.Range.Delete()
}

# ---------------------------------------------------------------
# Step 6: merge the "January 9" heading's first run with the trailing
# ", 2014" run into a single run "January 9, 2014" (matches the diff,
# which collapses those two runs into one).
# ---------------------------------------------------------------
$r = Find-Range("January 9, 2014")
$xml = "<w:p $wns><w:r><w:t>January 9, 2014</w:t></w:r></w:p>"
$r.InsertXML($xml)
